# Auto-generated Excel COM-interop edit script
# Applies profit/price recalculation updates across multiple job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 5
$ws.Range("H5").Value = 271
$ws.Range("I5").Value = 225.2
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 225.2
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -110.2
$ws.Range("N5").Value = -730

# Row 15
$ws.Range("H15").Value = 1134.2439
$ws.Range("I15").Value = 1134.2439
$ws.Range("K15").Value = 3402.7317
$ws.Range("M15").Value = -3233.7317

# Row 40
$ws.Range("H40").Value = 5294
$ws.Range("I40").Value = 6761.25
$ws.Range("K40").Value = 6761.25
$ws.Range("M40").Value = -6586.25

# Row 62
$ws.Range("H62").Value = 3901.0435
$ws.Range("I62").Value = 3342.805
$ws.Range("K62").Value = 3342.805
$ws.Range("M62").Value = -2718.805

# Row 65
$ws.Range("H65").Value = 3901.0435
$ws.Range("I65").Value = 3342.805
$ws.Range("K65").Value = 16714.025
$ws.Range("M65").Value = -13594.025

# Row 107
$ws.Range("H107").Value = 1479.1578
$ws.Range("I107").Value = 1065
$ws.Range("K107").Value = 1065
$ws.Range("M107").Value = 855

# Row 112
$ws.Range("H112").Value = 5000.0938
$ws.Range("J112").Value = 5341.552
$ws.Range("L112").Value = 16024.656
$ws.Range("N112").Value = -18240.656

# Row 131
$ws.Range("H131").Value = 3762.9473
$ws.Range("I131").Value = 2906.0625
$ws.Range("K131").Value = 8718.1875
$ws.Range("M131").Value = -3678.1875

# Row 132
$ws.Range("H132").Value = 2434134.8
$ws.Range("I132").Value = 2564473
$ws.Range("J132").Value = 1153.6666
$ws.Range("K132").Value = 7693419
$ws.Range("L132").Value = 3460.9998
$ws.Range("M132").Value = -7690889
$ws.Range("N132").Value = -8520.9998

# Row 135
$ws.Range("H135").Value = 1201.15
$ws.Range("I135").Value = 564
$ws.Range("K135").Value = 5076
$ws.Range("M135").Value = -2541

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 13334.214
$ws.Range("I32").Value = 14042.291
$ws.Range("K32").Value = 14042.291
$ws.Range("M32").Value = -13755.291

# Row 74
$ws.Range("H74").Value = 117956.42
$ws.Range("I74").Value = 125567.4
$ws.Range("J74").Value = 26624.75
$ws.Range("K74").Value = 125567.4
$ws.Range("L74").Value = 26624.75
$ws.Range("M74").Value = -124693.4
$ws.Range("N74").Value = -28372.75

# Row 77
$ws.Range("H77").Value = 117956.42
$ws.Range("I77").Value = 125567.4
$ws.Range("J77").Value = 26624.75
$ws.Range("K77").Value = 627837
$ws.Range("L77").Value = 133123.75
$ws.Range("M77").Value = -623469
$ws.Range("N77").Value = -141859.75

# Row 119
$ws.Range("H119").Value = 69805.836
$ws.Range("J119").Value = 69805.836
$ws.Range("L119").Value = 69805.836
$ws.Range("N119").Value = -79481.836

# Row 122
$ws.Range("H122").Value = 1353.5294
$ws.Range("I122").Value = 1360.8788
$ws.Range("K122").Value = 4082.6364
$ws.Range("M122").Value = -1632.6364

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 94
$ws.Range("H94").Value = 2223.0588
$ws.Range("I94").Value = 1374.4166
$ws.Range("J94").Value = 4259.8
$ws.Range("K94").Value = 1374.4166
$ws.Range("L94").Value = 4259.8
$ws.Range("M94").Value = -923.4166
$ws.Range("N94").Value = -5161.8

# Row 134
$ws.Range("H134").Value = 2452.303
$ws.Range("I134").Value = 1661.909
$ws.Range("K134").Value = 4985.727000000001
$ws.Range("M134").Value = -2450.727000000001

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 4764516
$ws.Range("I31").Value = 9092070
$ws.Range("J31").Value = 4206.4
$ws.Range("K31").Value = 9092070
$ws.Range("L31").Value = 4206.4
$ws.Range("M31").Value = -9091775
$ws.Range("N31").Value = -4796.4

# Row 34
$ws.Range("H34").Value = 4764516
$ws.Range("I34").Value = 9092070
$ws.Range("J34").Value = 4206.4
$ws.Range("K34").Value = 9092070
$ws.Range("L34").Value = 4206.4
$ws.Range("M34").Value = -9091868
$ws.Range("N34").Value = -4610.4

# Row 132
$ws.Range("H132").Value = 36747.25
$ws.Range("I132").Value = 42513.457
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 127540.371
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -125010.371
$ws.Range("N132").Value = -11510

# Row 135
$ws.Range("H135").Value = 119997.46
$ws.Range("J135").Value = 119997.46
$ws.Range("L135").Value = 119997.46
$ws.Range("N135").Value = -130137.46

# Row 138
$ws.Range("H138").Value = 119999
$ws.Range("J138").Value = 119999
$ws.Range("L138").Value = 119999
$ws.Range("N138").Value = -130279

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 6974.25
$ws.Range("I3").Value = 2632.8333
$ws.Range("J3").Value = 19998.5
$ws.Range("K3").Value = 7898.499899999999
$ws.Range("L3").Value = 59995.5
$ws.Range("M3").Value = -7786.499899999999
$ws.Range("N3").Value = -60219.5

# Row 12
$ws.Range("H12").Value = 406.25
$ws.Range("I12").Value = 525
$ws.Range("J12").Value = 366.66666
$ws.Range("K12").Value = 1575
$ws.Range("L12").Value = 1099.99998
$ws.Range("M12").Value = -1402
$ws.Range("N12").Value = -1445.99998

# Row 37
$ws.Range("H37").Value = 42220.375
$ws.Range("J37").Value = 42220.375
$ws.Range("L37").Value = 126661.125
$ws.Range("N37").Value = -126885.125

# Row 52
$ws.Range("H52").Value = 5127.5
$ws.Range("J52").Value = 5127.5
$ws.Range("L52").Value = 15382.5
$ws.Range("N52").Value = -15914.5

# Row 56
$ws.Range("H56").Value = 6598.4
$ws.Range("I56").Value = 6598.4
$ws.Range("K56").Value = 6598.4
$ws.Range("M56").Value = -6068.4

# Row 81
$ws.Range("H81").Value = 4699.936
$ws.Range("J81").Value = 4829.5454
$ws.Range("L81").Value = 14488.6362
$ws.Range("N81").Value = -16734.6362

# Row 84
$ws.Range("H84").Value = 4699.936
$ws.Range("J84").Value = 4829.5454
$ws.Range("L84").Value = 43465.9086
$ws.Range("N84").Value = -54697.9086

# Row 112
$ws.Range("H112").Value = 3662.25
$ws.Range("I112").Value = 2949.6667
$ws.Range("K112").Value = 8849.000100000001
$ws.Range("M112").Value = -7741.000100000001

# Row 131
$ws.Range("H131").Value = 119875.39
$ws.Range("I131").Value = 355025.5
$ws.Range("J131").Value = 2300.3333
$ws.Range("K131").Value = 1065076.5
$ws.Range("L131").Value = 6900.999899999999
$ws.Range("M131").Value = -1060036.5
$ws.Range("N131").Value = -16980.9999

# Row 132
$ws.Range("H132").Value = 1570.3572
$ws.Range("J132").Value = 937.25
$ws.Range("L132").Value = 8435.25
$ws.Range("N132").Value = -13495.25

# Row 139
$ws.Range("H139").Value = 1524.8572
$ws.Range("J139").Value = 1718.75
$ws.Range("L139").Value = 5156.25
$ws.Range("N139").Value = -15436.25

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 9
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 132
$ws.Range("H132").Value = 3182.75
$ws.Range("I132").Value = 3061.6
$ws.Range("K132").Value = 9184.799999999999
$ws.Range("M132").Value = -6654.799999999999

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 3007.8462
$ws.Range("I7").Value = 2972.5557
$ws.Range("K7").Value = 2972.5557
$ws.Range("M7").Value = -2860.5557

# Row 46
$ws.Range("H46").Value = 2948.913
$ws.Range("I46").Value = 664
$ws.Range("J46").Value = 5043.4165
$ws.Range("K46").Value = 664
$ws.Range("L46").Value = 5043.4165
$ws.Range("M46").Value = -476
$ws.Range("N46").Value = -5419.4165

# Row 126
$ws.Range("H126").Value = 3007.8462
$ws.Range("I126").Value = 2972.5557
$ws.Range("K126").Value = 8917.667099999999
$ws.Range("M126").Value = -6447.667099999999

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Range("H107").Value = 727.5833
$ws.Range("I107").Value = 538.3
$ws.Range("J107").Value = 1674
$ws.Range("K107").Value = 1614.9
$ws.Range("L107").Value = 5022
$ws.Range("M107").Value = 305.1000000000001
$ws.Range("N107").Value = -8862

# Row 113
$ws.Range("H113").Value = 665.2308
$ws.Range("I113").Value = 399.47058
$ws.Range("K113").Value = 1198.41174
$ws.Range("M113").Value = 971.58826

# Row 119
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
